# Apply the CodeSystem-patient-status metadata update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Bump the Version value (row 3) from 5.0.0 -> 6.0.0
$ws.Cells.Item(3, 2).Value2 = "6.0.0"

# 2. Update the Date value (row 8) to the new publication timestamp
$ws.Cells.Item(8, 2).Value2 = "2022-01-21T20:46:54+00:00"

# 3. Publisher (row 9) now carries an actual value instead of being blank
$ws.Cells.Item(9, 2).Value2 = "Alvearie Team"

# 4. Remove the duplicated "Contact / No display for ContactDetail" row
#    (originally rows 10 and 11 were identical); delete the second one so
#    everything below shifts up by one row.
$ws.Rows.Item(11).Delete()

# 5. The remaining "Contact" row (now row 10) becomes the new
#    "Jurisdiction / United States of America" row.
$ws.Cells.Item(10, 1).Value2 = "Jurisdiction"
$ws.Cells.Item(10, 2).Value2 = "United States of America"

# 6. Case Sensitive row (now row 14 after the deletion) gets a value of "true"
#    Assigning the literal string "true" directly gets auto-coerced to a
#    Boolean by Excel's smart cell-entry (same as typing TRUE into a cell),
#    so build it as a formula that yields the text "true" and then convert
#    that formula to its static value via copy / paste-values - this keeps
#    the result a genuine text cell instead of a boolean.
$caseSensitiveCell = $ws.Cells.Item(14, 2)
$caseSensitiveCell.Formula = '="tr"&"ue"'
$caseSensitiveCell.Copy() | Out-Null
$caseSensitiveCell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
